# Insert a new column "allele" before the existing "OR (CI = 95%)" column (currently column H).
# This shifts columns H:N -> I:O, and we then populate the new column H with header + values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at H; everything from H onward (including the
# "OR (CI = 95%)" header and its data) moves one column to the right (I..O).
$ws.Columns.Item(8).Insert()

# New header for the inserted column.
$ws.Range("H1").Value = "allele"

# Allele values for rows 2-28 (data rows), matching the source diff.
$alleles = @("A","G","C","T","G","C","C","A","A","G","A","A","T","A","T","C","T","G","A","T","T","A","G","C","T","T","A")

for ($i = 0; $i -lt $alleles.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $alleles[$i]
}
